# Refresh the live crypto-price snapshot on Sheet1 (rows 2-51) to match the
# latest scrape. Column D ("Price") holds numeric-looking text such as
# "104.30" or "1.00" -- if it is assigned as a plain string, Excel's COM
# layer auto-coerces it to a Number and silently drops the formatting
# (trailing zeros, "."-grouped thousands, etc.). Prefixing the literal with
# an apostrophe (a doubled '' inside a single-quoted PS string) forces Excel
# to keep storing it as Text, exactly like the source workbook; resetting
# Style back to 'Normal' afterwards clears the transient quote-prefix style
# so the cell's style index is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''47.907.08'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").Value = '''2.484.54'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -1.53%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''317.46'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").Value = '''104.30'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -4.96%  '

$ws.Range("E7").Value = '  -2.66%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -3.16%  '

$ws.Range("D10").Value = '''38.72'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -4.24%  '

$ws.Range("D11").Value = '''20.19'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -1.04%  '

$ws.Range("D12").Value = '''0.0798'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -3.00%  '

$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("D14").Value = '''7.03'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -3.35%  '

$ws.Range("D15").Value = '''2.875.68'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.49%  '

$ws.Range("D16").Value = '''2.496.03'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").Value = '''0.823'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -3.52%  '

$ws.Range("D18").Value = '''47.839.02'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -0.72%  '

$ws.Range("B19").Value = 'ImmutableX'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D19").Value = '''2.90'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +8.03%  '

$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = '''12.65'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -4.79%  '

$ws.Range("D21").Value = '''6.52'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.83%  '

$ws.Range("D22").Value = '''0.0₃0926'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -2.37%  '

$ws.Range("D23").Value = '''277.35'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("D24").Value = '''70.61'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -2.08%  '

$ws.Range("E25").Value = '  -3.65%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").Value = '''25.54'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("D28").Value = '''2.23'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -6.36%  '

$ws.Range("D29").Value = '''9.55'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -5.74%  '

$ws.Range("E30").Value = '  -4.06%  '

$ws.Range("E31").Value = '  -4.25%  '

$ws.Range("D32").Value = '''48.97'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -1.43%  '

$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").Value = '''18.97'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -4.00%  '

$ws.Range("D35").Value = '''5.23'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -2.99%  '

$ws.Range("D36").Value = '''0.0766'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -2.51%  '

$ws.Range("E37").Value = '  -2.94%  '

$ws.Range("D38").Value = '''4.48'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -4.66%  '

$ws.Range("E39").Value = '  -5.41%  '

$ws.Range("E40").Value = '  -1.54%  '

$ws.Range("E41").Value = '  -1.44%  '

$ws.Range("D42").Value = '''119.62'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -2.38%  '

$ws.Range("E43").Value = '  -3.01%  '

$ws.Range("E44").Value = '  -1.22%  '

$ws.Range("D45").Value = '''1.985.23'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.29%  '

$ws.Range("E46").Value = '  -1.74%  '

$ws.Range("E47").Value = '  +1.15%  '

$ws.Range("D48").Value = '''2.05'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("E49").Value = '  -2.63%  '

$ws.Range("D50").Value = '''5.07'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -2.67%  '

$ws.Range("D51").Value = '''78.67'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -1.82%  '
